# Atualização de bases das ligas, do dia: 15-06-2024 às 21:10
#
# The source data rows got re-sorted / re-matched against updated match ids.
# For each affected worksheet row, every column from B (id) through AD
# (PL_AhUnder) gets replaced with the values that used to live in a sibling
# row, while column A (the running row index) is left untouched.
#
# Rather than hard-coding every single new value (error prone), we capture
# the "before" values for B:AD of each row involved in a permutation group,
# then write them back out according to the row permutation describes by
# the diff (target row <- source row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (2) through AD (30) inclusive hold the per-match data that moves
# between rows; column A (1) is the stable row index and stays untouched.
$firstCol = 2
$lastCol = 30

function Get-RowValues($row) {
    $vals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value()
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c]
    }
}

# Each group lists the rows whose B:AD content gets permuted among
# themselves. $map[$targetRow] = $sourceRow means: after the edit,
# $targetRow's B:AD values equal what $sourceRow's B:AD values were
# *before* the edit.
$groups = @(
    @{ Rows = @(32,33);      Map = @{ 32 = 33; 33 = 32 } },
    @{ Rows = @(45,46);      Map = @{ 45 = 46; 46 = 45 } },
    @{ Rows = @(85,86);      Map = @{ 85 = 86; 86 = 85 } },
    @{ Rows = @(133,134);    Map = @{ 133 = 134; 134 = 133 } },
    @{ Rows = @(141,142);    Map = @{ 141 = 142; 142 = 141 } },
    @{ Rows = @(205,206);    Map = @{ 205 = 206; 206 = 205 } },
    @{ Rows = @(246,247);    Map = @{ 246 = 247; 247 = 246 } },
    @{ Rows = @(289,290);    Map = @{ 289 = 290; 290 = 289 } },
    @{ Rows = @(293,294,295); Map = @{ 293 = 295; 294 = 293; 295 = 294 } },
    @{ Rows = @(296,297);    Map = @{ 296 = 297; 297 = 296 } },
    @{ Rows = @(302,303,304); Map = @{ 302 = 304; 303 = 302; 304 = 303 } }
)

foreach ($group in $groups) {
    # Snapshot the original values of every row in this group before
    # writing anything back, since several rows feed into each other.
    $snapshot = @{}
    foreach ($r in $group.Rows) {
        $snapshot[$r] = Get-RowValues $r
    }

    foreach ($r in $group.Rows) {
        $srcRow = $group.Map[$r]
        Set-RowValues $r $snapshot[$srcRow]
    }
}
